$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "66.986.73"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.119.49"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("E4").Value = "  -0.04%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "577.32"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "173.02"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +2.45%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("E9").Value = "  -2.45%  "

$ws.Range("E10").Value = "  -1.11%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.483"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000247"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -1.26%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "37.09"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("E14").Value = "  -1.28%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.635.81"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "66.914.63"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("E17").Value = "  -0.37%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "3.117.23"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "16.28"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "476.47"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +2.20%  "

$ws.Range("E21").Value = "  -0.50%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "7.93"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +5.44%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "84.01"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "

$ws.Range("E24").Value = "  +2.76%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "2.29"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -3.17%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "10.11"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("E27").Value = "  +0.01%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "7.90"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -1.32%  "

$ws.Range("E29").Value = "  -1.44%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "2.68"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "28.55"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  +0.29%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0946"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -7.63%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  -0.61%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.974"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -3.06%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "47.12"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "50.19"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("E39").Value = "  -2.51%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.313"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -2.00%  "

$ws.Range("E41").Value = "  +0.67%  "

$ws.Range("E42").Value = "  -0.16%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "2.817.39"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +1.61%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "384.04"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("E45").Value = "  -1.87%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "2.55"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -10.01%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "135.48"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "24.87"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("E51").Value = "  -0.78%  "
